# Update target cluster expression data with new TPM values (rows 2-9)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.344207
$ws.Range("H2").Value = 1.032621
$ws.Range("I2").Value = 0.6985282229833164
$ws.Range("J2").Value = 0.6985282229833165
$ws.Range("M2").Value = 9.579981
$ws.Range("N2").Value = 28.739943
$ws.Range("O2").Value = 0.2856968652430823
$ws.Range("P2").Value = 0.2856968652430823
$ws.Range("Q2").Value = 3.297496520067
$ws.Range("R2").Value = 29.677468680603
$ws.Range("S2").Value = 0.1995673235901543
$ws.Range("T2").Value = 0.1995673235901543

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.344207
$ws.Range("H3").Value = 1.032621
$ws.Range("I3").Value = 0.6985282229833164
$ws.Range("J3").Value = 0.6985282229833165
$ws.Range("O3").Value = 0.2881547578255002
$ws.Range("P3").Value = 0.2881547578255002
$ws.Range("Q3").Value = 3.325865372592999
$ws.Range("R3").Value = 29.932788353337
$ws.Range("S3").Value = 0.2012842309280345
$ws.Range("T3").Value = 0.2012842309280345

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.344207
$ws.Range("H4").Value = 1.032621
$ws.Range("I4").Value = 0.6985282229833164
$ws.Range("J4").Value = 0.6985282229833165
$ws.Range("M4").Value = 9.479142000000001
$ws.Range("N4").Value = 28.437426
$ws.Range("O4").Value = 0.2826896164610391
$ws.Range("P4").Value = 0.2826896164610391
$ws.Range("Q4").Value = 3.262787030394
$ws.Range("R4").Value = 29.365083273546
$ws.Range("S4").Value = 0.1974666754423649
$ws.Range("T4").Value = 0.1974666754423649

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.344207
$ws.Range("H5").Value = 1.032621
$ws.Range("I5").Value = 0.6985282229833164
$ws.Range("J5").Value = 0.6985282229833165
$ws.Range("M5").Value = 4.810455999999999
$ws.Range("N5").Value = 14.431368
$ws.Range("O5").Value = 0.1434587604703784
$ws.Range("P5").Value = 0.1434587604703784
$ws.Range("Q5").Value = 1.655792628392
$ws.Range("R5").Value = 14.902133655528
$ws.Range("S5").Value = 0.1002099930227627
$ws.Range("T5").Value = 0.1002099930227627

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.1485533333333333
$ws.Range("H6").Value = 0.44566
$ws.Range("I6").Value = 0.3014717770166836
$ws.Range("J6").Value = 0.3014717770166836
$ws.Range("M6").Value = 9.579981
$ws.Range("N6").Value = 28.739943
$ws.Range("O6").Value = 0.2856968652430823
$ws.Range("P6").Value = 0.2856968652430823
$ws.Range("Q6").Value = 1.42313811082
$ws.Range("R6").Value = 12.80824299738
$ws.Range("S6").Value = 0.08612954165292801
$ws.Range("T6").Value = 0.08612954165292801

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.1485533333333333
$ws.Range("H7").Value = 0.44566
$ws.Range("I7").Value = 0.3014717770166836
$ws.Range("J7").Value = 0.3014717770166836
$ws.Range("O7").Value = 0.2881547578255002
$ws.Range("P7").Value = 0.2881547578255002
$ws.Range("Q7").Value = 1.435381579446667
$ws.Range("R7").Value = 12.91843421502
$ws.Range("S7").Value = 0.08687052689746565
$ws.Range("T7").Value = 0.08687052689746565

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.1485533333333333
$ws.Range("H8").Value = 0.44566
$ws.Range("I8").Value = 0.3014717770166836
$ws.Range("J8").Value = 0.3014717770166836
$ws.Range("M8").Value = 9.479142000000001
$ws.Range("N8").Value = 28.437426
$ws.Range("O8").Value = 0.2826896164610391
$ws.Range("P8").Value = 0.2826896164610391
$ws.Range("Q8").Value = 1.40815814124
$ws.Range("R8").Value = 12.67342327116
$ws.Range("S8").Value = 0.08522294101867418
$ws.Range("T8").Value = 0.08522294101867418

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.1485533333333333
$ws.Range("H9").Value = 0.44566
$ws.Range("I9").Value = 0.3014717770166836
$ws.Range("J9").Value = 0.3014717770166836
$ws.Range("M9").Value = 4.810455999999999
$ws.Range("N9").Value = 14.431368
$ws.Range("O9").Value = 0.1434587604703784
$ws.Range("P9").Value = 0.1434587604703784
$ws.Range("Q9").Value = 0.7146092736533333
$ws.Range("R9").Value = 6.431483462879999
$ws.Range("S9").Value = 0.04324876744761574
$ws.Range("T9").Value = 0.04324876744761576
